$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)

$rowData = @{
    4 = @(2, "Holden", 0.9955799824068393, 1.017732027400873, 0.9955436118146983, 1.017732027400873, 0.9955436118146983, 0.9882202176661367, 1.009744198633545, 0.9955693795253444, 0.9955436118146983, 0.9955799824068393, 1.006656004903856, 1.006656004903856, 1.007685402813752, 1.002951873874137, 1.002951873874137, 1.001099808359277, 1.001099808359277, 1.000398236241239)
    5 = @(3, "Rizzie Spiral", 1.009210368197925, 0.9631798984689328, 1.009195395637806, 0.9631798984689328, 1.009195395637806, 1.024563853974221, 0.9797454824664167, 1.009205996931287, 1.009195395637806, 1.009210368197925, 0.9861951333334289, 0.9861951333334289, 0.9840452497110915, 0.993861887434888, 0.993861887434888, 0.9976952644856175, 0.9976952644856175, 0.9991834992794315)
    6 = @(2, "RotRing OmegaMax-90", 1.002663928160471, 0.9893454911472179, 1.00266308129584, 0.9893454911472179, 1.00266308129584, 1.00710397695084, 0.994139824433027, 1.002663682799783, 1.00266308129584, 1.002663928160471, 0.9960047096538444, 0.9960047096538444, 0.9953830812469052, 0.998224166867843, 0.998224166867843, 0.9993338954748423, 0.9993338954748423, 0.9997633307978631)
    7 = @(3, "Equal Angle", 0.9949976620749289, 1.020058886930838, 0.9949629863832838, 1.020058886930838, 0.9949629863832838, 0.9866670408645514, 1.011024362716139, 0.9949875531700276, 0.9949629863832838, 0.9949976620749289, 1.007528274502884, 1.007528274502884, 1.008693637240635, 1.003339845129684, 1.003339845129684, 1.001245630443084, 1.001245630443084, 1.000449748689962)
    8 = @(4, "Tilt Rotate", 0.9827026386248929, 1.069436575668115, 0.9825296455986897, 1.069436575668115, 0.9825296455986897, 0.953906656312446, 1.038150064817806, 0.982652201239355, 0.9825296455986897, 0.9827026386248929, 1.026069607146504, 1.026069607146504, 1.030096426370271, 1.011556286630566, 1.011556286630566, 1.004299626372597, 1.004299626372597, 1.001562963710217)
    9 = @(5, "CLR", 0.9995066336133799, 1.001965169718086, 0.9995124412818406, 1.001965169718086, 0.9995124412818406, 0.9986832532268719, 1.001082188924882, 0.9995083272848116, 0.9995124412818406, 0.9995066336133799, 1.000735901665733, 1.000735901665733, 1.000851330752116, 1.000328081537769, 1.000328081537769, 1.000124171473787, 1.000124171473787, 1.000043002341645)
    10 = @(6, "Rizzie Hex", 0.9999565484448533, 1.000172314683491, 0.9999575928523742, 1.000172314683491, 0.9999575928523742, 0.999883933753003, 1.000095015443496, 0.9999568539160034, 0.9999575928523742, 0.9999565484448533, 1.000064431564172, 1.000064431564172, 1.000074626190613, 1.00002881866024, 1.00002881866024, 1.000011012208273, 1.000011012208273, 1.00000370984887)
    11 = @(7, "Matthies Hex", 0.9990445177703934, 1.003811004089877, 0.9990521657258106, 1.003811004089877, 0.9990521657258106, 0.9974505937605196, 1.002097821896416, 0.9990467475825771, 0.9990521657258106, 0.9990445177703934, 1.001427760930135, 1.001427760930135, 1.001651114585562, 1.000635895862027, 1.000635895862027, 1.000239963327973, 1.000239963327973, 1.000083808470932)
    12 = @(8, "Tilt Rotate_Partial", 0.9824455286893193, 1.070472737844497, 0.9822671298941084, 1.070472737844497, 0.9822671298941084, 0.9532220602718665, 1.038718702374367, 0.9823935152468015, 0.9822671298941084, 0.9824455286893193, 1.026459133266908, 1.026459133266908, 1.030545656302728, 1.011728465475975, 1.011728465475975, 1.004363131580508, 1.004363131580508, 1.001586612386827)
    13 = @(9, "RotRing OmegaMax-60", 1.001065495688795, 0.9956609743520584, 1.001119424644171, 0.9956609743520584, 1.001119424644171, 1.002831052301122, 0.9976260171683459, 1.001081219007493, 1.001119424644171, 1.001065495688795, 0.9983632350204268, 0.9983632350204268, 0.9981174957363997, 0.9992819648950082, 0.9992819648950082, 0.9997413298322988, 0.9997413298322988, 0.9998973638603309)
    14 = @(10, "Equal Angle_Partial", 0.9948276541052654, 1.020770161894738, 0.9947711098315786, 1.020770161894738, 0.9947711098315786, 0.9862178508421059, 1.01141050234737, 0.9948111689263149, 0.9947711098315786, 0.9948276541052654, 1.007798908000002, 1.007798908000002, 1.009002772782458, 1.003456308610527, 1.003456308610527, 1.00128500891579, 1.00128500891579, 1.000468074657896)
    15 = @(11, "Rizzie Hex_Partial", 1.003459787529534, 0.9861540377158861, 1.003464559894419, 0.9861540377158861, 1.003464559894419, 1.009225194939853, 0.9923858246697717, 1.00346117914602, 1.003464559894419, 1.003459787529534, 0.9948069126227099, 0.9948069126227099, 0.9939998833050638, 0.9976927950466129, 0.9976927950466131, 0.9991357362585647, 0.9991357362585647, 0.9996917639825806)
    16 = @(12, "ND Single", 0.9696726700000003, 1.1217596, 0.9693574699999995, 1.1217596, 0.9693574699999995, 0.9191871600000003, 1.0668948, 0.9695807699999998, 0.9693574699999995, 0.9696726700000003, 1.045716135, 1.045716135, 1.05277569, 1.020263246666667, 1.020263246666667, 1.0075368025, 1.0075368025, 1.002742078333333)
    17 = @(13, "RD Single", 1.0113912, 0.9526597999999999, 1.012634, 0.9526597999999999, 1.012634, 1.0301398, 0.9742505199999999, 1.0117535, 1.012634, 1.0113912, 0.9820255, 0.9820255, 0.97943384, 0.9922283333333333, 0.9922283333333333, 0.99732975, 0.99732975, 0.9988048033333333)
    18 = @(14, "TD Single", 1.0193616, 0.92240399, 1.0194663, 0.92240399, 1.0194663, 1.051611, 0.9573464299999999, 1.0193921, 1.0194663, 1.0193616, 0.9708827950000001, 0.9708827950000001, 0.9663706733333334, 0.9870772966666667, 0.9870772966666665, 0.9951745474999999, 0.9951745474999999, 0.9982635699999999)
    19 = @(15, "Morris Single", 1.0218585, 0.9130504, 1.0215194, 0.9130504, 1.0215194, 1.0583539, 0.95209925, 1.0217596, 1.0215194, 1.0218585, 0.96745445, 0.96745445, 0.9623360499999999, 0.9854761, 0.9854761000000001, 0.9944869250000001, 0.9944869250000001, 0.9981068416666666)
    20 = @(16, "Ring Perpendicular to ND", 0.9888877309589039, 1.044285950410959, 0.989001912328767, 1.044285950410959, 0.989001912328767, 0.9703455435616437, 1.024383695068493, 0.9889210199999997, 0.989001912328767, 0.9888877309589039, 1.016586840684931, 1.016586840684931, 1.019185792146118, 1.00739186456621, 1.00739186456621, 1.002794376506849, 1.002794376506849, 1.000970975388128)
    21 = @(17, "Ring Perpendicular to RD", 1.003952777368421, 0.9847224236842105, 1.003579299473684, 0.9847224236842105, 1.003579299473684, 1.010611881052632, 0.9915108994736841, 1.003843893157895, 1.003579299473684, 1.003952777368421, 0.9943376005263158, 0.9943376005263158, 0.9933953668421053, 0.9974181668421053, 0.9974181668421053, 0.9989584500000001, 0.9989584500000001, 0.9997035290350879)
    22 = @(18, "Ring Perpendicular to TD", 1.006396526842105, 0.9744232884210525, 1.006389938421053, 0.9744232884210525, 1.006389938421053, 1.017058677894737, 0.9859312868421052, 1.006394599473684, 1.006389938421053, 1.006396526842105, 0.9904099076315789, 0.9904099076315789, 0.9889170340350876, 0.9957365845614033, 0.9957365845614033, 0.9983999230263156, 0.9983999230263156, 0.9994323863157893)
    23 = @(19, "OffsetFTD", 1.002724931667661, 0.9897079552393377, 1.00229955870053, 0.9897079552393377, 1.00229955870053, 1.007347508454681, 0.9942409311711703, 1.002600916710594, 1.00229955870053, 1.002724931667661, 0.9962164434534996, 0.9962164434534996, 0.9955579393593897, 0.9982441485358432, 0.9982441485358432, 0.9992580010770149, 0.9992580010770149, 0.9998203003239957)
    24 = @(20, "OffsetATD", 1.000606737637396, 0.9973229015443735, 1.00078182513898, 0.9973229015443735, 1.00078182513898, 1.001584642071828, 0.9985681129572284, 1.00065779050995, 1.00078182513898, 1.000606737637396, 0.9989648195908847, 0.9989648195908847, 0.9988325840463327, 0.9995704881069165, 0.9995704881069164, 0.9998733223649323, 0.9998733223649323, 0.9999203349766259)
    25 = @(21, "OffsetF45", 1.001628178870439, 0.9930606403511493, 1.001926814300016, 0.9930606403511493, 1.001926814300016, 1.004284944662976, 0.9962524615799423, 1.001715254269505, 1.001926814300016, 1.001628178870439, 0.997344409610794, 0.997344409610794, 0.9969804269338435, 0.9988718778405344, 0.9988718778405344, 0.9996356119554047, 0.9996356119554047, 0.9998113823390044)
    26 = @(22, "OffsetA45", 1.000463792942587, 0.9983240710954566, 1.000338328286865, 0.9983240710954566, 1.000338328286865, 1.001260677753924, 0.9990492061484643, 1.000427214557603, 1.000338328286865, 1.000463792942587, 0.999393932019022, 0.999393932019022, 0.9992790233955028, 0.9997087307749696, 0.9997087307749696, 0.9998661301529435, 0.9998661301529435, 0.9999772151308166)
    27 = @(23, "OffsetFRD", 1.009885808848352, 0.9602868156588219, 1.01000477125186, 0.9602868156588219, 1.01000477125186, 1.026339519567837, 0.978185285694646, 1.009920491676321, 1.01000477125186, 1.009885808848352, 0.9850863122535871, 0.9850863122535871, 0.9827859700672734, 0.9933924652530114, 0.9933924652530114, 0.9975455417527236, 0.9975455417527236, 0.9991037821163063)
    28 = @(24, "OffsetARD", 1.002610240226877, 0.9896280724699763, 1.002561924092682, 0.9896280724699763, 1.002561924092682, 1.006969844814552, 0.9942842547364354, 1.00259615183624, 1.002561924092682, 1.002610240226877, 0.9961191563484266, 0.9961191563484266, 0.9955075224777629, 0.9982667455965117, 0.9982667455965117, 0.9993405402205543, 0.9993405402205543, 0.9997750813627938)
    29 = @(25, "Gaussian Quadrature", 0.9981862414686992, 1.00724474927493, 0.9981934321675275, 1.00724474927493, 0.9981934321675275, 0.9951619627560849, 1.003986286561895, 0.9981883482093247, 0.9981934321675275, 0.9981862414686992, 1.002715495371814, 1.002715495371814, 1.003139092435175, 1.001208140970385, 1.001208140970385, 1.000454463769671, 1.000454463769671, 1.000160170073077)
    30 = @(26, "Michael-CCHex", 1.002114922908906, 0.9916039621993677, 1.002070398310387, 0.9916039621993677, 1.002070398310387, 1.005648245832334, 0.9953718739847448, 1.002101925236739, 1.002070398310387, 1.002114922908906, 0.9968594425541368, 0.9968594425541368, 0.9963635863643395, 0.9985964278062203, 0.9985964278062203, 0.999464920432262, 0.999464920432262, 0.9998185547454131)
    31 = @(27, "Michael-SNHex", 1.007580556009474, 0.9696859848469154, 1.007574876011742, 0.9696859848469154, 1.007574876011742, 1.020215843011948, 0.9833259797293258, 1.007578876826781, 1.007574876011742, 1.007580556009474, 0.9886332704281948, 0.9886332704281948, 0.9868641735285718, 0.9949471389560438, 0.9949471389560438, 0.9981040732199682, 0.9981040732199682, 0.999327019406031)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item([int]$r, 1).Value = $vals[0]
    $ws.Cells.Item([int]$r, 2).Value = $vals[1]
    for ($i = 2; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, $i + 1).Value = $vals[$i]
    }
}